$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Junio de 2020 a las 04:05"

# --- Swap Guatemala/Honduras rows (71/72) with updated case numbers ---
$ws.Range("A71").Value = "Honduras"
$ws.Range("B71").Value = 5362
$ws.Range("C71").Value = 160
$ws.Range("D71").Value = 549
$ws.Range("E71").Value = 4596
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 5
$ws.Range("H71").Value = 217

$ws.Range("A72").Value = "Guatemala"
$ws.Range("B72").Value = 5336
$ws.Range("C72").Value = 249
$ws.Range("D72").Value = 795
$ws.Range("E72").Value = 4425
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 8
$ws.Range("H72").Value = 116

# --- Swap Santa Lucia/Belice rows (201/202) with updated case numbers ---
$ws.Range("A201").Value = "Belice"
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 16
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 2

$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("B202").Value = 18
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 18
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

# --- Update Corea del Sur row (53) numbers ---
$ws.Range("B53").Value = 11541
$ws.Range("C53").Value = 38
$ws.Range("D53").Value = 10446
$ws.Range("E53").Value = 823
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 272

# --- Update Haiti row (91) numbers ---
$ws.Range("B91").Value = 2226
$ws.Range("C91").Value = 102
$ws.Range("D91").Value = 24
$ws.Range("E91").Value = 2157
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 45
